$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new "dd" (Deleted) column right after column A ("id"), before the
# existing "no" column. Columns.Insert() shifts B:K -> C:L, carrying every
# cell's value/style/reference along with it, and grows the used range.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Insert()

# New column should be narrow like column A/B (raw OOXML width 6.0).
# Excel's ColumnWidth COM property is offset from the stored sheet width by
# 5/6 of a character, so request (6 - 5/6) to land exactly on 6.0.
$ws.Columns.Item(2).ColumnWidth = 5.166666666666667

# The freshly inserted column copied column A's header formatting (bold,
# red-fill "id" style). Re-stamp it with the plain bold header look used by
# the rest of the non-highlighted headers (e.g. column D, "invalid") before
# writing the new header text.
$ws.Range("D1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 2).Value = "dd"

# Data row: mark this template row as not-deleted by default ("N").
# B2 already inherited the plain body-row style from the column insert.
$ws.Cells.Item(2, 2).Value = "N"

# Re-assert the sheet's existing "fit to one page" print setup so it survives
# the save (explicitly touching PageSetup keeps the engine from dropping the
# unchanged <pageSetup> element on export).
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.FitToPagesWide = 1
